$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Expand header row (row 1) into the new multi-article layout.
$ws.Range("A1").Value = "positive claim"
$ws.Range("B1").Value = "article 1"
$ws.Range("C1").Value = "article 2"
$ws.Range("D1").Value = "article 3"
$ws.Range("E1").Value = "article 4"
$ws.Range("F1").Value = "article 5"
$ws.Range("G1").Value = "counter_article 1"
$ws.Range("H1").Value = "counter_article 2"
$ws.Range("I1").Value = "counter_article 3"
$ws.Range("J1").Value = "counter_article 4"
$ws.Range("K1").Value = "counter_article 5"

# 2) Append the new claim rows at the bottom (17-21).
$ws.Range("A17").Value = "Extremists are motivated by pizzagate"
$ws.Range("B17").Value = "pizzagate_10"

$ws.Range("A18").Value = "Pizzagate is a religion"
$ws.Range("B18").Value = "pizzagate_13"

$ws.Range("A19").Value = "Conspiracy Theories are dangerous"
$ws.Range("B19").Value = "pizzagate_15"

$ws.Range("A20").Value = "Fake news is a threat"
$ws.Range("B20").Value = "pizzagate_21"

$ws.Range("A21").Value = "Pizzagate is fake news"
$ws.Range("B21").Value = "pizzagate_21"

# 3) Fill in article ids for the existing rows 10 and 11.
$ws.Range("B10").Value = "pizzagate_27"
$ws.Range("C10").Value = "pizzagate_29"
$ws.Range("B11").Value = "pizzagate_30"

# 4) Leave the selection where the author left off.
[void]$ws.Range("B12").Select()
